$wb = $excel.ActiveWorkbook

# The UniProt-derived "EGF-like N" domain labels are renamed to the more
# specific "I-EGF N" nomenclature used by the integrin beta chain (ITGB1 /
# ITGB3) transcripts; the underlying start/end coordinates were also
# re-derived after switching to "import Utils.X" (qualified) instead of
# "import X" (unqualified), which changed a couple of resolved offsets.

$sheetNames = @("ENST00000302278", "ENST00000559488")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    for ($i = 1; $i -le 4; $i++) {
        $row = 6 + $i
        $ws.Cells.Item($row, 9).Value = "I-EGF $i"
    }
}

$ws2 = $wb.Worksheets.Item("ENST00000302278")
$ws2.Cells.Item(7, 7).Value = 466
$ws2.Cells.Item(10, 8).Value = 631

$ws4 = $wb.Worksheets.Item("ENST00000559488")
$ws4.Cells.Item(7, 7).Value = 463
$ws4.Cells.Item(10, 8).Value = 625
